# Apply updates to the ActivePlayers_NZL sheet (all_formats_raw)
# Columns: A=NAME, B=PLAYER_ID, C=TEST, D=ODI, E=T20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_formats_raw")

# Row 5 - T A Blundell: D5 2 -> 3
$ws.Range("D5").Value = 3

# Row 7 - C J Bowes: D7 1 -> 2, E7 (empty) -> 2
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2

# Row 10 - M S Chapman: E10 44 -> 46
$ws.Range("E10").Value = 46

# Row 17 - M J Henry: D17 66 -> 67
$ws.Range("D17").Value = 67

# Row 20 - T W M Latham: D20 124 -> 125, E20 18 -> 20
$ws.Range("D20").Value = 125
$ws.Range("E20").Value = 20

# Row 21 - B G Lister: E21 1 -> 3
$ws.Range("E21").Value = 3

# Row 23 - A F Milne: E23 35 -> 37
$ws.Range("E23").Value = 37

# Row 24 - D J Mitchell: D24 20 -> 21, E24 44 -> 46
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 46

# Row 25 - J D S Neesham: E25 60 -> 62
$ws.Range("E25").Value = 62

# Row 26 - H M Nicholls: D26 61 -> 62
$ws.Range("D26").Value = 62

# Row 29 - R Ravindra: D29 1 -> 2, E29 6 -> 8
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 8

# Row 33 - T L Seifert: E33 40 -> 42
$ws.Range("E33").Value = 42

# Row 34 - H B Shipley: D34 4 -> 5, E34 (empty) -> 2
$ws.Range("D34").Value = 5
$ws.Range("E34").Value = 2

# Row 35 - I S Sodhi: D35 40 -> 41, E35 91 -> 93
$ws.Range("D35").Value = 41
$ws.Range("E35").Value = 93

# Row 39 - B M Tickner: D39 10 -> 11
$ws.Range("D39").Value = 11

# Row 43 - W A Young: D43 9 -> 10
$ws.Range("D43").Value = 10
